$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $ws.Range("D$r").Value = 44232.51400987701
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Range("D$r").Value = 44232.49295149306
}

for ($r = 30; $r -le 37; $r++) {
    $ws.Range("D$r").Value = 44232.47190263889
}
